$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "canonical SMILES" column (D)
$ws.Range("D2").Value = 'canonical SMILES'
$ws.Range("D3").Value = 'COc1cccc(c1)[NH+]=c2c3ccccc3[nH]cn2'
$ws.Range("D4").Value = 'COc1cccc(c1)N=c2c3ccccc3[nH]cn2'
$ws.Range("D5").Value = 'COc1cccc(c1)Nc2c3ccccc3ncn2'
$ws.Range("D6").Value = 'COc1cccc(c1)Nc2c3ccccc3nc[nH+]2'
$ws.Range("D7").Value = 'COc1cccc(c1)[NH2+]c2c3ccccc3ncn2'
$ws.Range("D8").Value = 'COc1cccc(c1)[NH2+]c2c3ccccc3[nH+]cn2'
$ws.Range("D9").Value = 'COc1cccc(c1)[N-]c2c3ccccc3ncn2'
$ws.Range("D10").Value = 'COc1cccc(c1)N=c2c3ccccc3nc[nH]2'
$ws.Range("D11").Value = 'COc1cccc(c1)Nc2c3ccccc3[nH+]c[nH+]2'
$ws.Range("D12").Value = 'COc1cccc(c1)[NH2+]c2c3ccccc3nc[nH+]2'
$ws.Range("D13").Value = 'COc1cccc(c1)[NH2+]c2c3ccccc3[nH+]c[nH+]2'

# Set column D width to match source width (closest achievable in this engine)
$ws.Columns.Item(4).ColumnWidth = 36
